$wb = $excel.ActiveWorkbook

$oldGuid = "3b4a67ad-5fcc-48cf-bd33-a0f76c157783"
$newGuid = "adcad4ac-7f22-493a-a9d3-a80be848f884"
$oldHash = "67c98bd21d4981b5901d9e82a2dc45ec9f57b8f1"
$newHash = "55b9b357dff3e4f0bb816b41b86038331490476c"

$newMdName = "$newGuid.md"
$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"

# Overview sheet: A2 = handoff md file name, D2 = latest handoff date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = "2016-44-18 10:44:41"
foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Column -eq 1) {
        $hl.TextToDisplay = $newMdName
    }
}

# zh-cn sheet: A2 = md file name, D2 = handoff xlf file name, E2 = handoff datetime
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhXlf
$wsZh.Range("E2").Value = "2016-03-18 10:44:38"
foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.Range.Column -eq 1) {
        $hl.TextToDisplay = $newMdName
    } elseif ($hl.Range.Column -eq 4) {
        $hl.TextToDisplay = $newZhXlf
    }
}

# de-de sheet: A2 = md file name, D2 = handoff xlf file name, E2 = handoff datetime
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeXlf
$wsDe.Range("E2").Value = "2016-03-18 10:44:41"
foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.Range.Column -eq 1) {
        $hl.TextToDisplay = $newMdName
    } elseif ($hl.Range.Column -eq 4) {
        $hl.TextToDisplay = $newDeXlf
    }
}
